$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Pagos" (F) and "Inscrições homologadas" (H) values for several rows

$ws.Range("F4").Value = 9
$ws.Range("H4").Value = 12

$ws.Range("F5").Value = 13
$ws.Range("H5").Value = 17

$ws.Range("F7").Value = 17
$ws.Range("H7").Value = 18

$ws.Range("F11").Value = 12
$ws.Range("H11").Value = 13

$ws.Range("F12").Value = 11
$ws.Range("H12").Value = 13

$ws.Range("F14").Value = 18
$ws.Range("H14").Value = 20
